$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record arrived for "Feria Lagunitas de Puerto Montt" /
# Albahaca. It belongs right under the header, so insert a fresh row at 182
# and push the existing rows 182-188 down to 183-189 (dimension grows to
# A1:R189).
$ws.Rows.Item(182).Insert()

# Populate the newly inserted row 182 with the new week's data.
$ws.Range("A182").Value = 4
$ws.Range("B182").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C182").Value = "Los Lagos"
$ws.Range("D182").Value = 45041
$ws.Range("E182").Value = 10
$ws.Range("F182").Value = 100112052
$ws.Range("G182").Value = "Albahaca"
$ws.Range("H182").Value = "Sin especificar"
$ws.Range("I182").Value = "Primera"
$ws.Range("J182").Value = 80
$ws.Range("K182").Value = 8000
$ws.Range("L182").Value = 8000
$ws.Range("M182").Value = 8000
$ws.Range("N182").Value = "$/docena de matas"
$ws.Range("O182").Value = "Región Metropolitana"
$ws.Range("P182").Value = 1333
$ws.Range("Q182").Value = 6
$ws.Range("R182").Value = "Hortaliza"
